$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 173, pushing existing rows 173:205 down to 174:206.
$ws.Rows(173).Insert()

# Populate the newly inserted row 173 with the new weekly record
$ws.Cells.Item(173, 1).Value = 4
$ws.Cells.Item(173, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(173, 3).Value = "Los Lagos"
$ws.Cells.Item(173, 4).Value = 44785
$ws.Cells.Item(173, 5).Value = 10
$ws.Cells.Item(173, 6).Value = 100112009
$ws.Cells.Item(173, 7).Value = "Acelga"
$ws.Cells.Item(173, 8).Value = "Sin especificar"
$ws.Cells.Item(173, 9).Value = "Primera"
$ws.Cells.Item(173, 10).Value = 200
$ws.Cells.Item(173, 11).Value = 1500
$ws.Cells.Item(173, 12).Value = 1500
$ws.Cells.Item(173, 13).Value = 1500
$ws.Cells.Item(173, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(173, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(173, 16).Value = 1000
$ws.Cells.Item(173, 17).Value = 1.5
$ws.Cells.Item(173, 18).Value = "Hortaliza"
